$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and title to reflect new "through" date
$ws.Name = "Through 2022-10-23"

# Update the label for the October row
$ws.Range("A11").Value = "October (through 10-23)"

# Update October row (row 11) values
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 39
$ws.Range("D11").Value = 50
$ws.Range("E11").Value = 52
$ws.Range("F11").Value = 39
$ws.Range("G11").Value = 111
$ws.Range("H11").Value = 147
$ws.Range("I11").Value = 85

# Update Total row (row 12) values
$ws.Range("B12").Value = 246
$ws.Range("C12").Value = 468
$ws.Range("D12").Value = 677
$ws.Range("E12").Value = 600
$ws.Range("F12").Value = 461
$ws.Range("G12").Value = 1012
$ws.Range("H12").Value = 1394
$ws.Range("I12").Value = 1362
